# hr_boot: adding the rockitBerlin platform
#
# The "RockITdigital" sheet describes how the job post should be
# forwarded to the rockIT digital platform. This edit:
#   - clears the "Vertrag" (contract) cell, it no longer applies
#   - switches the "&" separators used when copy-pasting the allowed
#     "Berufserfahrung" / "Kategorie(n)" values to ","  (rockIT's own
#     listing doesn't accept "+"/"&" as a separator)
#   - widens column A and G so the new values/labels are easier to read
#   - leaves the user on the RockITdigital tab, scrolled/selected on the
#     "Bewerbung auf folgender Website" block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RockITdigital")

# Make this the active sheet/tab.
$ws.Activate()

# D3 "Vertrag" no longer has a value.
$ws.Range("D3").Value = ""

# E3 "Berufserfahrung" / F3 "Kategorie(n)": swap "+"/"&" for ","
$ws.Range("E3").Value = "Keine , Erste"
$ws.Range("F3").Value = "Marketing/ Online Marketing , Suchmaschinenmarketing (SEM/ SEA)"

# Widen column A and column G to fit the (new) content better.
$ws.Range("A1").EntireColumn.ColumnWidth = 38 - 0.8333333333
$ws.Range("G1").EntireColumn.ColumnWidth = 75 - 0.8333333333

# Leave the selection on the "Bewerbung auf folgender Website" header.
$ws.Range("G1:H1").Select()
